# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that are plain numeric-looking strings (e.g. '555.50')
# are assigned with a leading apostrophe so Excel keeps them as TEXT
# instead of auto-converting them to numbers (matching the source data,
# which stores all Price/Volume values as text).

$ws.Range('D2').Value = '63.025.75'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.682.70'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('D5').Value = "'" + '555.50'
$ws.Range('E5').Value = '  -2.76%  '
$ws.Range('D6').Value = "'" + '158.91'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'" + '0.593'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('E11').Value = '  -4.17%  '
$ws.Range('D12').Value = "'" + '5.38'
$ws.Range('E12').Value = '  -7.16%  '
$ws.Range('D13').Value = '3.156.93'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').Value = "'" + '26.37'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '62.898.53'
$ws.Range('E15').Value = '  -1.52%  '
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('D17').Value = '2.682.69'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = "'" + '11.92'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('D19').Value = "'" + '4.63'
$ws.Range('E19').Value = '  -3.72%  '
$ws.Range('D20').Value = "'" + '345.65'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').Value = "'" + '6.30'
$ws.Range('E21').Value = '  -4.71%  '
$ws.Range('D22').Value = "'" + '1.00'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = "'" + '0.510'
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('D24').Value = "'" + '63.36'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range('D26').Value = "'" + '0.999'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').Value = "'" + '8.22'
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('D28').Value = '0.0₃0860'
$ws.Range('E28').Value = '  -6.20%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = "'" + '7.29'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = "'" + '1.37'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = "'" + '165.50'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'" + '1.49'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'" + '4.89'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = "'" + '19.55'
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('D37').Value = "'" + '1.79'
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').Value = "'" + '348.95'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = "'" + '6.36'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = "'" + '0.958'
$ws.Range('E40').Value = '  -3.31%  '
$ws.Range('D41').Value = "'" + '4.02'
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('D42').Value = "'" + '38.29'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').Value = "'" + '20.43'
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('D44').Value = "'" + '20.88'
$ws.Range('E44').Value = '  -5.06%  '
$ws.Range('E45').Value = '  -1.30%  '
$ws.Range('D46').Value = "'" + '0.0563'
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = "'" + '11.03'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('E50').Value = '  -3.41%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'" + '129.12'
$ws.Range('E51').Value = '  -3.80%  '
